$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values for weeks 12 and 13
$ws.Range("B13").Value = 377
$ws.Range("B14").Value = 450

# Add new rows for weeks 14-17
$data = @(
    @(14, 456),
    @(15, 354),
    @(16, 257),
    @(17, 3)
)

$row = 15
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
